$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64 (was Cortinarius dalecarlicus / Siljansspindling) -> becomes Ramaria rufescens / Fjällfotad fingersvamp
$ws.Range("A64").Value = 111998588
$ws.Range("B64").Value = 89103
$ws.Range("E64").Value = 233196
$ws.Range("F64").Value = "Fjällfotad fingersvamp"
$ws.Range("G64").Value = "Ramaria rufescens"
$ws.Range("H64").Value = "(Schaeff.) Corner"
$ws.Range("Q64").Value = 485479
$ws.Range("R64").Value = 6995888
$ws.Range("Z64").Value = "16:31"
$ws.Range("AB64").Value = "16:31"

# Row 65 (was Ramaria rufescens / Fjällfotad fingersvamp) -> becomes Cortinarius dalecarlicus / Siljansspindling
$ws.Range("A65").Value = 111998587
$ws.Range("B65").Value = 85327
$ws.Range("E65").Value = 3595
$ws.Range("F65").Value = "Siljansspindling"
$ws.Range("G65").Value = "Cortinarius dalecarlicus"
$ws.Range("H65").Value = "Brandrud"
$ws.Range("Q65").Value = 485439
$ws.Range("R65").Value = 6995893
$ws.Range("Z65").Value = "16:38"
$ws.Range("AB65").Value = "16:38"

# Row 66: only the Taxonsorteringsordning (B) value changes
$ws.Range("B66").Value = 90817

# Row 67 (was Boletopsis leucomelaena / Grangråticka) -> becomes Ramaria pallida / Blek fingersvamp
$ws.Range("A67").Value = 111998589
$ws.Range("B67").Value = 89098
$ws.Range("D67").Value = "NT"
$ws.Range("E67").Value = 256756
$ws.Range("F67").Value = "Blek fingersvamp"
$ws.Range("G67").Value = "Ramaria pallida"
$ws.Range("H67").Value = "(Schaeff.) Ricken"
$ws.Range("Q67").Value = 485479
$ws.Range("Z67").Value = "16:20"
$ws.Range("AB67").Value = "16:20"

# Row 68 (was Ramaria pallida / Blek fingersvamp) -> becomes Boletopsis leucomelaena / Grangråticka
$ws.Range("A68").Value = 111998584
$ws.Range("B68").Value = 90803
$ws.Range("D68").Value = "VU"
$ws.Range("E68").Value = 150
$ws.Range("F68").Value = "Grangråticka"
$ws.Range("G68").Value = "Boletopsis leucomelaena"
$ws.Range("H68").Value = "(Pers.) Fayod"
$ws.Range("Q68").Value = 485433
$ws.Range("Z68").Value = "16:53"
$ws.Range("AB68").Value = "16:53"

# Row 69: only the Taxonsorteringsordning (B) value changes
$ws.Range("B69").Value = 84955
